$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MOSIP-14369: the "is_active" column (E2:E5) was stored as a TRUE() boolean
# formula; fix it so the cells hold the literal text value "TRUE" instead.
foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Formula = '="TRUE"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Reflect the edited range in the sheet's active selection.
$ws.Range("E2:E5").Select()
